$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so numeric-looking
# strings (e.g. "583.47") are written back as text, matching the
# original inlineStr cell type instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns for rows with changed values
$ws.Range("D2").Value = "62.400.19"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.451.90"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "583.47"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").Value = "143.54"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "2.447.27"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "5.22"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "0.346"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").Value = "26.49"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "0.0000177"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "2.896.35"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "62.255.38"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "2.446.96"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "10.80"
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D20").Value = "7.14"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "327.58"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "1.96"
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "9.21"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "589.82"
$ws.Range("E27").Value = "  -5.37%  "
$ws.Range("D28").Value = "0.0₃0975"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "8.02"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "1.90"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").Value = "0.135"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").Value = "4.94"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "1.45"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "0.378"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "154.08"
$ws.Range("E39").Value = "  +5.35%  "
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").Value = "142.96"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").Value = "3.65"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "0.0₆0254"
$ws.Range("E48").Value = "  +14.83%  "
$ws.Range("D49").Value = "0.608"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "0.0525"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "19.90"
$ws.Range("E51").Value = "  -2.08%  "

# Row 40/41 swap: RenderToken <-> EthereumClassic (with refreshed values)
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.45"
$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "5.31"
$ws.Range("E41").Value = "  +0.96%  "

# Row 42/43 swap: Stacks <-> OKB (with refreshed values)
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "42.88"
$ws.Range("E42").Value = "  +1.76%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.72"
$ws.Range("E43").Value = "  -0.94%  "

# Restore column D back to its original (default/General) style so no
# stray number-format is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
